$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three more scraped match rows were appended under the existing ones.
# They duplicate the data already present in rows 4, 3 and 2 (in that
# order), so copy/paste those rows instead of retyping the values - this
# keeps the numeric-looking text (e.g. "0", "15", "187.50") stored as
# text, exactly like the existing rows, instead of being re-interpreted
# as numbers.
$xlPasteAll = -4104

$ws.Range("A4:K4").Copy()
$ws.Range("A5:K5").PasteSpecial($xlPasteAll)

$ws.Range("A3:K3").Copy()
$ws.Range("A6:K6").PasteSpecial($xlPasteAll)

$ws.Range("A2:K2").Copy()
$ws.Range("A7:K7").PasteSpecial($xlPasteAll)

$excel.CutCopyMode = 0
